$d = $word.ActiveDocument

# Find the start of the "Proposal B" heading ("Proposal B:") - this is
# where the whole second proposal section begins, right after the page
# break that separates it from Proposal A.
$searchRange = $d.Content.Duplicate
$found = $searchRange.Find.Execute("Proposal B:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $deleteRange = $d.Range($searchRange.Start, $d.Content.End)
    $deleteRange.Delete()
}
